$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "IPC"
$ws.Range("D4").Value = "Monthly CPI"
$ws.Range("F4").Value = "Could not fetch the inhambane CPI, the url is http://www.ine.gov.mz/estatisticas/estatisticas-economicas/indice-de-preco-no-consumidor/quadros/inhambane"
$ws.Range("G4").Value = "2022-09-08 19:44:11"
